$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$emuPerPt = 12700.0

# Helper: fully replace the (single-row/single-column) text of a table cell
# without leaving stray leftover runs -- TextRange.Text on a cell that already
# contains multiple runs only overwrites the first run in this host, so we
# add a fresh row, drop the old row (taking its runs with it) and then set
# the text on the new, run-less paragraph that is left behind.
function Set-SoleCellText($table, [string]$text) {
    $table.Rows.Add() | Out-Null
    $table.Rows(1).Delete() | Out-Null
    $table.Cell(1,1).Shape.TextFrame.TextRange.Text = $text
}

# ---------------------------------------------------------------------------
# 1) Duplicate the existing "Table 12" (em/ab diagram table) twice BEFORE
#    touching its text, so the two brand-new tables inherit its run/paragraph
#    formatting (lang, dirty, table style, grid/row extension ids, ...).
# ---------------------------------------------------------------------------
$tbl12Shape = $s.Shapes.Item(2)

$newShape1 = $tbl12Shape.Duplicate().Item(1)
$newShape1.Name = "Table 16"
$newShape1.Left = 473240 / $emuPerPt
$newShape1.Top = 3926064 / $emuPerPt
$newShape1.Width = 1833356 / $emuPerPt
$newShape1.Height = 410363 / $emuPerPt
Set-SoleCellText $newShape1.Table "em0:EventManager"

$newShape2 = $tbl12Shape.Duplicate().Item(1)
$newShape2.Name = "Table 17"
$newShape2.Left = 2472474 / $emuPerPt
$newShape2.Top = 3929598 / $emuPerPt
$newShape2.Width = 1833356 / $emuPerPt
$newShape2.Height = 410363 / $emuPerPt
Set-SoleCellText $newShape2.Table "em1:EventManager"

# ---------------------------------------------------------------------------
# 2) Update "Table 12" itself: ab0:AddressBook -> em0:EventManager
# ---------------------------------------------------------------------------
Set-SoleCellText $tbl12Shape.Table "em0:EventManager"

# ---------------------------------------------------------------------------
# 3) Remove the old "Table 11" / "Table 14" (ab1:/ab0:AddressBook) shapes.
# ---------------------------------------------------------------------------
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Table 11" -or $sh.Name -eq "Table 14") {
        $sh.Delete()
    }
}
